# Scheduled runner update: refresh market-board derived price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the crafting-class
# sheets. Some HQ profit cells (M/N) appear or disappear depending on
# whether an HQ price is available for the item on that pass, so those
# cells are added or cleared accordingly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1095.25
$ws.Range("I12").Value = 861.4
$ws.Range("K12").Value = 861.4
$ws.Range("M12").Value = -691.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 5899.8335
$ws.Range("I34").Value = 5899.8335
$ws.Range("K34").Value = 5899.8335
$ws.Range("M34").Value = -5696.8335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 5899.8335
$ws.Range("I36").Value = 5899.8335
$ws.Range("K36").Value = 5899.8335
$ws.Range("M36").Value = -5184.8335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3197.6
$ws.Range("I113").Value = 2997
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2997
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 257
$ws.Range("N113").Value = -10508

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6149.533
$ws.Range("J116").Value = 8373.6
$ws.Range("L116").Value = 8373.6
$ws.Range("N116").Value = -15257.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1371.9166
$ws.Range("I132").Value = 1371.9166
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4115.7498
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1585.7498
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1060.2667
$ws.Range("I135").Value = 674.7692
$ws.Range("J135").Value = 3566
$ws.Range("K135").Value = 6072.922799999999
$ws.Range("L135").Value = 32094
$ws.Range("M135").Value = -3537.922799999999
$ws.Range("N135").Value = -37164

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4070.3958
$ws.Range("I138").Value = 1992
$ws.Range("J138").Value = 4312.07
$ws.Range("K138").Value = 5976
$ws.Range("L138").Value = 12936.21
$ws.Range("M138").Value = -836
$ws.Range("N138").Value = -23216.21

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1392.5
$ws.Range("I61").Value = 1392.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1392.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1180.5
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1429.6666
$ws.Range("I132").Value = 1360.381
$ws.Range("K132").Value = 4081.143
$ws.Range("M132").Value = -1551.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1392.5
$ws.Range("I136").Value = 1392.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4177.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1627.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 299980
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 299980
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 299980
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -309780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2024
$ws.Range("I134").Value = 2488.6667
$ws.Range("J134").Value = 630
$ws.Range("K134").Value = 7466.000100000001
$ws.Range("L134").Value = 1890
$ws.Range("M134").Value = -4931.000100000001
$ws.Range("N134").Value = -6960

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 59713.285
$ws.Range("J62").Value = 82398.60000000001
$ws.Range("L62").Value = 82398.60000000001
$ws.Range("N62").Value = -83646.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 59713.285
$ws.Range("J65").Value = 82398.60000000001
$ws.Range("L65").Value = 411993
$ws.Range("N65").Value = -418233

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 16897.842
$ws.Range("I99").Value = 15158.4
$ws.Range("J99").Value = 18830.555
$ws.Range("K99").Value = 15158.4
$ws.Range("L99").Value = 18830.555
$ws.Range("M99").Value = -13660.4
$ws.Range("N99").Value = -21826.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 16897.842
$ws.Range("I126").Value = 15158.4
$ws.Range("J126").Value = 18830.555
$ws.Range("K126").Value = 45475.2
$ws.Range("L126").Value = 56491.665
$ws.Range("M126").Value = -43005.2
$ws.Range("N126").Value = -61431.665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 606.6842
$ws.Range("I107").Value = 646.6667
$ws.Range("J107").Value = 603.25714
$ws.Range("K107").Value = 1940.0001
$ws.Range("L107").Value = 1809.77142
$ws.Range("M107").Value = -20.00009999999997
$ws.Range("N107").Value = -5649.77142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1706
$ws.Range("I132").Value = 1082
$ws.Range("K132").Value = 9738
$ws.Range("M132").Value = -7208

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6852.846
$ws.Range("I80").Value = 6498
$ws.Range("J80").Value = 7266.8335
$ws.Range("K80").Value = 6498
$ws.Range("L80").Value = 7266.8335
$ws.Range("M80").Value = -5500
$ws.Range("N80").Value = -9262.833500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 6852.846
$ws.Range("I83").Value = 6498
$ws.Range("J83").Value = 7266.8335
$ws.Range("K83").Value = 32490
$ws.Range("L83").Value = 36334.1675
$ws.Range("M83").Value = -27498
$ws.Range("N83").Value = -46318.1675

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4999.857
$ws.Range("I126").Value = 4999
$ws.Range("K126").Value = 14997
$ws.Range("M126").Value = -12527

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1489.5714
$ws.Range("I7").Value = 1574.4
$ws.Range("J7").Value = 1277.5
$ws.Range("K7").Value = 1574.4
$ws.Range("L7").Value = 1277.5
$ws.Range("M7").Value = -1462.4
$ws.Range("N7").Value = -1501.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1915.5
$ws.Range("I40").Value = 1915.5
$ws.Range("K40").Value = 1915.5
$ws.Range("M40").Value = -1779.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9749.75
$ws.Range("I122").Value = 9749.75
$ws.Range("K122").Value = 29249.25
$ws.Range("M122").Value = -26799.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1489.5714
$ws.Range("I126").Value = 1574.4
$ws.Range("J126").Value = 1277.5
$ws.Range("K126").Value = 4723.200000000001
$ws.Range("L126").Value = 3832.5
$ws.Range("M126").Value = -2253.200000000001
$ws.Range("N126").Value = -8772.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8100
$ws.Range("J62").Value = 8625
$ws.Range("L62").Value = 8625
$ws.Range("N62").Value = -9873

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 8100
$ws.Range("J65").Value = 8625
$ws.Range("L65").Value = 43125
$ws.Range("N65").Value = -49365

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5079.75
$ws.Range("I122").Value = 5079.75
$ws.Range("K122").Value = 15239.25
$ws.Range("M122").Value = -12789.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1797.0588
$ws.Range("I132").Value = 1603.5385
$ws.Range("K132").Value = 4810.6155
$ws.Range("M132").Value = -2280.6155
